# Add a new "2022-Q3" sheet (becomes the 2nd tab, right after the "总计"
# summary sheet) and push the existing quarter sheets down, and add a
# corresponding summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Match the look & feel (outline + page margins) used by the other
# quarter sheets.
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# ---------------------------------------------------------------------
# 3. Data rows (2022-Q3 fund holdings).
# columns: A idx, B code, C name, D scale, E stock-pos, F pos-pct, G value, H rank
# ---------------------------------------------------------------------
$rows = @(
    ,@(0,  "000727", "融通健康产业灵活配置混合A",       "22.64", "93.67", "5.00", "1.1320", 7)
    ,@(1,  "009274", "融通健康产业灵活配置混合C",       "17.64", "93.67", "5.00", "0.8820", 7)
    ,@(2,  "010709", "安信医药健康主题股票A",           "16.29", "94.32", "3.89", "0.6337", 8)
    ,@(3,  "006218", "富国生物医药科技混合A",           "7.10",  "89.32", "8.16", "0.5794", 1)
    ,@(4,  "010710", "安信医药健康主题股票C",           "13.91", "94.32", "3.89", "0.5411", 8)
    ,@(5,  "100016", "富国天源沪港深平衡混合A",         "4.99",  "70.11", "7.96", "0.3972", 1)
    ,@(6,  "011308", "富国生物医药科技混合C",           "1.52",  "89.32", "8.16", "0.1240", 1)
    ,@(7,  "014867", "摩根士丹利华鑫优悦安和混合C",     "1.32",  "93.41", "5.41", "0.0714", 10)
    ,@(8,  "011404", "融通鑫新成长混合C",               "1.75",  "94.07", "4.06", "0.0710", 6)
    ,@(9,  "009893", "摩根士丹利华鑫优悦安和混合A",     "0.87",  "93.41", "5.41", "0.0471", 10)
    ,@(10, "014220", "恒越医疗健康精选混合A",           "0.72",  "88.76", "4.94", "0.0356", 5)
    ,@(11, "009246", "大摩ESG量化混合",                 "2.69",  "84.29", "1.01", "0.0272", 6)
    ,@(12, "011403", "融通鑫新成长混合A",               "0.39",  "94.07", "4.06", "0.0158", 6)
    ,@(13, "014221", "恒越医疗健康精选混合C",           "0.29",  "88.76", "4.94", "0.0143", 5)
    ,@(14, "014931", "富国天源沪港深平衡混合C",         "0.00",  "70.11", "7.96", $null,    1)
)

# Text-typed columns (B fund code, D,E,F,G) need the "@" number format
# pre-applied so that values like "000727" / "5.00" survive as literal
# text instead of being coerced to numbers (losing leading zeros /
# trailing zeros).  Apply it to the whole block up front.
$q3.Range("B2:B16").NumberFormat = "@"
$q3.Range("D2:G16").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    if ($row[6] -eq $null) {
        $q3.Cells.Item($r, 7).NumberFormat = "General"
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        $q3.Cells.Item($r, 7).Value = $row[6]
    }
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# The "@" text format leaves a style index behind even though the
# target file has no styling on these cells at all; strip it back down
# to the default "no style" by pasting the format from a pristine,
# never-touched cell on the same sheet.
$q3.Range("Z1").Copy()
$q3.Range("B2:B16").PasteSpecial(-4122)
$q3.Range("D2:G16").PasteSpecial(-4122)
$q3.Range("A1").Select()

# Bold/bordered style used for the header row and the leading index
# column — reuse the exact style already present on the "总计" sheet so
# the same style id is produced.
$summary.Range("A2").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2:A16").PasteSpecial(-4122)
$q3.Range("A2:A16").PasteSpecial(-4122)

# Re-apply values after the format paste (PasteSpecial with formats only
# should not disturb them, but re-assert to be safe for the header text
# values which must remain the original strings, not copied from A2).
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# ---------------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: insert the 2022-Q3 row at the
#    top of the data block and push the rest down by one row.
# ---------------------------------------------------------------------
$summaryRows = @(
    ,@("2022-Q3", 15, 4.57)
    ,@("2022-Q2", 16, 5.91)
    ,@("2022-Q1", 6,  3.18)
    ,@("2021-Q4", 5,  0.84)
    ,@("2021-Q3", 1,  1.04)
)

$r = 2
$idx = 0
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $idx
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
    $r++
    $idx++
}

# A6 is a brand new cell — give it the same bold/bordered style as the
# rest of column A on this sheet.
$summary.Range("A2").Copy()
$summary.Range("A6").PasteSpecial(-4122)

$summary.Activate()
$summary.Range("A1").Select()

Write-Output "edit complete"
